$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 111112660
$ws.Range("I62").Value = 142858270
$ws.Range("J62").Value = 3006
$ws.Range("K62").Value = 142858270
$ws.Range("L62").Value = 3006
$ws.Range("M62").Value = -142857646

$ws.Range("H64").Value = 3732.359
$ws.Range("I64").Value = 3524.75
$ws.Range("J64").Value = 3756.0857
$ws.Range("K64").Value = 3524.75
$ws.Range("L64").Value = 3756.0857
$ws.Range("M64").Value = -3276.75
$ws.Range("N64").Value = -4252.0857

$ws.Range("H65").Value = 111112660
$ws.Range("I65").Value = 142858270
$ws.Range("J65").Value = 3006
$ws.Range("K65").Value = 714291350
$ws.Range("L65").Value = 15030
$ws.Range("M65").Value = -714288230

$ws.Range("H67").Value = 3732.359
$ws.Range("I67").Value = 3524.75
$ws.Range("J67").Value = 3756.0857
$ws.Range("K67").Value = 3524.75
$ws.Range("L67").Value = 3756.0857
$ws.Range("M67").Value = -2666.75
$ws.Range("N67").Value = -5472.0857

$ws.Range("H68").Value = 29294.25
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 29294.25
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 29294.25
$ws.Range("N68").Value = -30792.25
$ws.Range("M68").ClearContents()

$ws.Range("H71").Value = 29294.25
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 29294.25
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 87882.75
$ws.Range("N71").Value = -95370.75
$ws.Range("M71").ClearContents()

$ws.Range("H112").Value = 3894.348
$ws.Range("I112").Value = 697.5
$ws.Range("J112").Value = 4198.8096
$ws.Range("K112").Value = 2092.5
$ws.Range("L112").Value = 12596.4288
$ws.Range("M112").Value = -984.5
$ws.Range("N112").Value = -14812.4288

$ws.Range("H141").Value = 15661.667
$ws.Range("I141").Value = 20410
$ws.Range("J141").Value = 9726.25
$ws.Range("K141").Value = 61230
$ws.Range("L141").Value = 29178.75
$ws.Range("M141").Value = -56050
$ws.Range("N141").Value = -39538.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 383.0435
$ws.Range("I2").Value = 380.75
$ws.Range("J2").Value = 398.33334
$ws.Range("K2").Value = 380.75
$ws.Range("L2").Value = 398.33334
$ws.Range("M2").Value = -267.75
$ws.Range("N2").Value = -624.33334

$ws.Range("H74").Value = 1273.0435
$ws.Range("I74").Value = 804
$ws.Range("J74").Value = 1343.4
$ws.Range("K74").Value = 804
$ws.Range("L74").Value = 1343.4
$ws.Range("M74").Value = 70
$ws.Range("N74").Value = -3091.4

$ws.Range("H77").Value = 1273.0435
$ws.Range("I77").Value = 804
$ws.Range("J77").Value = 1343.4
$ws.Range("K77").Value = 4020
$ws.Range("L77").Value = 6717
$ws.Range("M77").Value = 348
$ws.Range("N77").Value = -15453

$ws.Range("H116").Value = 383.0435
$ws.Range("I116").Value = 380.75
$ws.Range("J116").Value = 398.33334
$ws.Range("K116").Value = 380.75
$ws.Range("L116").Value = 398.33334
$ws.Range("M116").Value = 1913.25
$ws.Range("N116").Value = -4986.33334

$ws.Range("H132").Value = 2719570.2
$ws.Range("I132").Value = 5436238
$ws.Range("J132").Value = 2902.4348
$ws.Range("K132").Value = 16308714
$ws.Range("L132").Value = 8707.304400000001
$ws.Range("M132").Value = -16306184
$ws.Range("N132").Value = -13767.3044

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 383.0435
$ws.Range("I3").Value = 380.75
$ws.Range("J3").Value = 398.33334
$ws.Range("K3").Value = 380.75
$ws.Range("L3").Value = 398.33334
$ws.Range("M3").Value = -266.75
$ws.Range("N3").Value = -626.33334

$ws.Range("H86").Value = 7144697.5
$ws.Range("I86").Value = 16668409
$ws.Range("J86").Value = 1913.25
$ws.Range("K86").Value = 16668409
$ws.Range("L86").Value = 1913.25
$ws.Range("M86").Value = -16667286
$ws.Range("N86").Value = -4159.25

$ws.Range("H89").Value = 7144697.5
$ws.Range("I89").Value = 16668409
$ws.Range("J89").Value = 1913.25
$ws.Range("K89").Value = 83342045
$ws.Range("L89").Value = 9566.25
$ws.Range("M89").Value = -83336429
$ws.Range("N89").Value = -20798.25

$ws.Range("H134").Value = 5299890
$ws.Range("I134").Value = 6182840
$ws.Range("J134").Value = 2190.4443
$ws.Range("K134").Value = 18548520
$ws.Range("L134").Value = 6571.3329
$ws.Range("M134").Value = -18545985

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1492.5834
$ws.Range("I16").Value = 1473.7273
$ws.Range("J16").Value = 1700
$ws.Range("K16").Value = 1473.7273
$ws.Range("L16").Value = 1700
$ws.Range("M16").Value = -1186.7273
$ws.Range("N16").Value = -2274

$ws.Range("H22").Value = 573.5454999999999
$ws.Range("I22").Value = 641.8333
$ws.Range("J22").Value = 491.6
$ws.Range("K22").Value = 641.8333
$ws.Range("L22").Value = 491.6
$ws.Range("M22").Value = -291.8333
$ws.Range("N22").Value = -1191.6

$ws.Range("H62").Value = 2623.75
$ws.Range("I62").Value = 2260
$ws.Range("J62").Value = 2987.5
$ws.Range("K62").Value = 2260
$ws.Range("L62").Value = 2987.5
$ws.Range("M62").Value = -1636
$ws.Range("N62").Value = -4235.5

$ws.Range("H65").Value = 2623.75
$ws.Range("I65").Value = 2260
$ws.Range("J65").Value = 2987.5
$ws.Range("K65").Value = 11300
$ws.Range("L65").Value = 14937.5
$ws.Range("M65").Value = -8180
$ws.Range("N65").Value = -21177.5

$ws.Range("H113").Value = 1492.5834
$ws.Range("I113").Value = 1473.7273
$ws.Range("J113").Value = 1700
$ws.Range("K113").Value = 1473.7273
$ws.Range("L113").Value = 1700
$ws.Range("M113").Value = 696.2727
$ws.Range("N113").Value = -6040

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 10213.454
$ws.Range("I68").Value = 399.5
$ws.Range("J68").Value = 12394.333
$ws.Range("K68").Value = 1198.5
$ws.Range("L68").Value = 37182.999
$ws.Range("M68").Value = -387.5
$ws.Range("N68").Value = -38804.999

$ws.Range("H71").Value = 10213.454
$ws.Range("I71").Value = 399.5
$ws.Range("J71").Value = 12394.333
$ws.Range("K71").Value = 3595.5
$ws.Range("L71").Value = 111548.997
$ws.Range("M71").Value = 460.5
$ws.Range("N71").Value = -119660.997
